$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $exactText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $pp = $doc.Paragraphs.Item($i)
        if ($pp.Range.Text -eq $exactText) {
            return $i
        }
    }
    throw "Could not find paragraph with text: $exactText"
}

# Note on step 5 below: inserting two chunks of text back-to-back via
# Range.InsertAfter on the same range coalesces them into a single <w:r>
# (same run formatting). To keep them as two separate runs (as in the
# target XML) we instead put the second chunk in a new paragraph right
# after the first, then delete the paragraph mark between them -- merging
# the two paragraphs back into one while keeping the two runs distinct
# (mirrors how real Word preserves run boundaries when a paragraph break is
# removed).

# 1. "What is an optimal sequential cut-off ..." question -> the following
#    empty paragraph receives the answer text.
$idx1 = Get-ParaIndexByText $d "What is an optimal sequential cut-off for both parallel algorithms? (Note that the optimal sequential cut-off can vary based on dataset size.)`r"
$target1 = $d.Paragraphs.Item($idx1 + 1)
$target1.Range.InsertAfter("The optimal sequential cut-off is based on the size of the image. So, for example, if the image is 1920x1080 then the optimal sequential cut-off will be somewhere around this area. If the sequential cut-off is any lower than the image size, then the program will crash.")

# 2. "For what range of data set sizes/ filter sizes ..." question -> the
#    following empty paragraph receives the answer text.
$idx2 = Get-ParaIndexByText $d "For what range of data set sizes/ filter sizes do your parallel programs perform well?`r"
$target2 = $d.Paragraphs.Item($idx2 + 1)
$target2.Range.InsertAfter("Based on the speed-up obtained above, the parallel programs performed excellently within the range of 3 – 21.")

# 3. "What is the maximum speedup obtainable ..." question -> the following
#    empty paragraph is split into three paragraphs, each holding one
#    sentence of the answer.
$idx3 = Get-ParaIndexByText $d "What is the maximum speedup obtainable with each parallel algorithm? How do they differ and why? How close is the speedup to the ideal expected?`r"
$target3 = $d.Paragraphs.Item($idx3 + 1)
$r3 = $target3.Range
$r3.InsertAfter("The maximum speed-up obtained was 68.7 for the mean filter, using the Girl.jpg and a window size of 17.")
$r3.InsertParagraphAfter()

$target3b = $d.Paragraphs.Item($idx3 + 2)
$target3b.Range.InsertAfter("The maximum speed-up obtained was 9.13 for the median filter, using the Girl.jpg and a window size of 17.")
$target3b.Range.InsertParagraphAfter()

$target3c = $d.Paragraphs.Item($idx3 + 3)
$target3c.Range.InsertAfter("The program seems to perform well on this image and window size for both filter.")

# 4. "How reliable are your measurements? ..." question -> the following
#    empty paragraph receives the answer text.
$idx4 = Get-ParaIndexByText $d "How reliable are your measurements? Are they any anomalies and can you explain why they occur?`r"
$target4 = $d.Paragraphs.Item($idx4 + 1)
$target4.Range.InsertAfter("It is difficult to determine where exactly to start and end timing calls, so there could be a slight inaccuracy in timing between the serial and parallel solutions. Furthermore, the programs that run in the background of the computer can hinder the performance of the algorithm, resulting in anomalies.  ")

# 5. The final (empty) paragraph of the document, right after "Conclusions",
#    receives two sentences as two separate runs.
$idx5 = Get-ParaIndexByText $d "Conclusions`r"
$target5 = $d.Paragraphs.Item($idx5 + 1)
$r5 = $target5.Range
$r5.InsertAfter("For both the median and mean filters, there was a speed-up when comparing the serial and parallel solutions. Therefore, ")
$r5.InsertParagraphAfter()

$target5b = $d.Paragraphs.Item($idx5 + 2)
$target5b.Range.InsertAfter("for the purpose of mean and median filtering in Java, it is worth designing a parallel solution. ")

# Merge the two paragraphs created in step 5 back into a single paragraph by
# deleting the paragraph mark between them, so the two sentences end up as
# two distinct runs inside one paragraph (matching the target XML).
$markRange = $d.Range($target5.Range.End - 1, $target5.Range.End)
$markRange.Delete()
